# Refresh the cryptocurrency price/volume snapshot (GitHub Actions scheduled update).
# A new coin (WrappedliquidstakedEther2.0) entered the top-50 ranking, shifting every
# row below it down by one and dropping the previous last row (NEARProtocol) off the
# bottom of the list; every other row keeps its coin/link but refreshes Price/Volume.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '''30.375.99'
$ws.Range("E2").Value = '''  +1.27%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '''1.926.67'
$ws.Range("E3").Value = '''  +0.92%  '

# Row 4: TetherUSD
$ws.Range("D4").Value = '''0.9998'
$ws.Range("E4").Value = '''  -0.05%  '

# Row 5: XRP
$ws.Range("D5").Value = '''0.8096'
$ws.Range("E5").Value = '''  +2.36%  '

# Row 6: BNB
$ws.Range("D6").Value = '''244.79'
$ws.Range("E6").Value = '''  +1.31%  '

# Row 7: USDC
$ws.Range("D7").Value = '''0.9997'
$ws.Range("E7").Value = '''  +0.00%  '

# Row 8: Cardano
$ws.Range("E8").Value = '''  +3.21%  '

# Row 9: Solana
$ws.Range("D9").Value = '''27.33'
$ws.Range("E9").Value = '''  +3.79%  '

# Row 10: Dogecoin
$ws.Range("D10").Value = '''0.07309'
$ws.Range("E10").Value = '''  +6.14%  '

# Row 11: Polygon
$ws.Range("D11").Value = '''0.7993'
$ws.Range("E11").Value = '''  +7.70%  '

# Row 12: TRON
$ws.Range("D12").Value = '''0.08095'
$ws.Range("E12").Value = '''  +1.19%  '

# Row 13: WrappedEther
$ws.Range("D13").Value = '''1.929.64'
$ws.Range("E13").Value = '''  +1.09%  '

# Row 14: Polkadot
$ws.Range("D14").Value = '''5.430'
$ws.Range("E14").Value = '''  +4.52%  '

# Row 15: Litecoin
$ws.Range("D15").Value = '''94.66'
$ws.Range("E15").Value = '''  +1.74%  '

# Row 16: WrappedBTC
$ws.Range("D16").Value = '''30.376.34'
$ws.Range("E16").Value = '''  +1.25%  '

# Row 17: Avalanche
$ws.Range("E17").Value = '''  +3.38%  '

# Row 18: Uniswap
$ws.Range("D18").Value = '''6.133'
$ws.Range("E18").Value = '''  +4.51%  '

# Row 19: BitcoinCash
$ws.Range("D19").Value = '''253.42'
$ws.Range("E19").Value = '''  +3.17%  '

# Row 20: ShibaInu
$ws.Range("D20").Value = '''0.000007899'
$ws.Range("E20").Value = '''  +2.12%  '

# Row 21: WrappedliquidstakedEther2.0
$ws.Range("B21").Value = '''WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = '''https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '''2.177.24'
$ws.Range("E21").Value = '''  +0.82%  '

# Row 22: Dai
$ws.Range("B22").Value = '''Dai'
$ws.Range("C22").Value = '''https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").Value = '''1.0000'
$ws.Range("E22").Value = '''  +0.02%  '

# Row 23: Chainlink
$ws.Range("B23").Value = '''Chainlink'
$ws.Range("C23").Value = '''https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D23").Value = '''8.140'
$ws.Range("E23").Value = '''  +19.11%  '

# Row 24: BinanceUSD
$ws.Range("B24").Value = '''BinanceUSD'
$ws.Range("C24").Value = '''https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D24").Value = '''0.9999'
$ws.Range("E24").Value = '''  -0.04%  '

# Row 25: Stellar
$ws.Range("B25").Value = '''Stellar'
$ws.Range("C25").Value = '''https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D25").Value = '''0.1634'
$ws.Range("E25").Value = '''  +17.58%  '

# Row 26: Cosmos
$ws.Range("B26").Value = '''Cosmos'
$ws.Range("C26").Value = '''https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '''9.583'
$ws.Range("E26").Value = '''  +3.93%  '

# Row 27: Monero
$ws.Range("B27").Value = '''Monero'
$ws.Range("C27").Value = '''https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '''167.79'
$ws.Range("E27").Value = '''  -0.32%  '

# Row 28: EthereumClassic
$ws.Range("B28").Value = '''EthereumClassic'
$ws.Range("C28").Value = '''https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '''19.16'
$ws.Range("E28").Value = '''  +1.48%  '

# Row 29: LidoDAOToken
$ws.Range("B29").Value = '''LidoDAOToken'
$ws.Range("C29").Value = '''https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").Value = '''2.163'
$ws.Range("E29").Value = '''  +6.34%  '

# Row 30: Toncoin
$ws.Range("B30").Value = '''Toncoin'
$ws.Range("C30").Value = '''https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = '''1.376'
$ws.Range("E30").Value = '''  +0.79%  '

# Row 31: PancakeSwap
$ws.Range("B31").Value = '''PancakeSwap'
$ws.Range("C31").Value = '''https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = '''1.548'
$ws.Range("E31").Value = '''  +2.06%  '

# Row 32: Filecoin
$ws.Range("B32").Value = '''Filecoin'
$ws.Range("C32").Value = '''https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '''4.359'
$ws.Range("E32").Value = '''  +1.09%  '

# Row 33: InternetComputer(DFINITY)
$ws.Range("B33").Value = '''InternetComputer(DFINITY)'
$ws.Range("C33").Value = '''https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").Value = '''4.162'
$ws.Range("E33").Value = '''  +2.01%  '

# Row 34: Hedera
$ws.Range("B34").Value = '''Hedera'
$ws.Range("C34").Value = '''https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '''0.05640'
$ws.Range("E34").Value = '''  +2.39%  '

# Row 35: ARBITRUM
$ws.Range("B35").Value = '''ARBITRUM'
$ws.Range("C35").Value = '''https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '''1.307'
$ws.Range("E35").Value = '''  +4.17%  '

# Row 36: ImmutableX
$ws.Range("B36").Value = '''ImmutableX'
$ws.Range("C36").Value = '''https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '''0.7464'
$ws.Range("E36").Value = '''  +1.74%  '

# Row 37: Frax
$ws.Range("B37").Value = '''Frax'
$ws.Range("C37").Value = '''https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D37").Value = '''1.004'
$ws.Range("E37").Value = '''  +0.48%  '

# Row 38: HuobiToken
$ws.Range("B38").Value = '''HuobiToken'
$ws.Range("C38").Value = '''https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D38").Value = '''2.719'
$ws.Range("E38").Value = '''  -0.05%  '

# Row 39: VeChain
$ws.Range("B39").Value = '''VeChain'
$ws.Range("C39").Value = '''https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '''0.01961'
$ws.Range("E39").Value = '''  +1.84%  '

# Row 40: MXToken
$ws.Range("B40").Value = '''MXToken'
$ws.Range("C40").Value = '''https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '''2.817'
$ws.Range("E40").Value = '''  +1.04%  '

# Row 41: TheSandbox
$ws.Range("B41").Value = '''TheSandbox'
$ws.Range("C41").Value = '''https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '''0.4516'
$ws.Range("E41").Value = '''  +2.22%  '

# Row 42: Aave
$ws.Range("B42").Value = '''Aave'
$ws.Range("C42").Value = '''https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '''74.33'
$ws.Range("E42").Value = '''  +2.85%  '

# Row 43: FraxShare
$ws.Range("B43").Value = '''FraxShare'
$ws.Range("C43").Value = '''https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '''6.004'
$ws.Range("E43").Value = '''  -2.19%  '

# Row 44: RenderToken
$ws.Range("B44").Value = '''RenderToken'
$ws.Range("C44").Value = '''https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '''1.942'
$ws.Range("E44").Value = '''  +3.73%  '

# Row 45: TrustWalletToken
$ws.Range("B45").Value = '''TrustWalletToken'
$ws.Range("C45").Value = '''https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").Value = '''0.8561'
$ws.Range("E45").Value = '''  +2.26%  '

# Row 46: PaxDollar
$ws.Range("B46").Value = '''PaxDollar'
$ws.Range("C46").Value = '''https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = '''0.9995'
$ws.Range("E46").Value = '''  -0.06%  '

# Row 47: Quant
$ws.Range("B47").Value = '''Quant'
$ws.Range("C47").Value = '''https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").Value = '''103.81'
$ws.Range("E47").Value = '''  +3.31%  '

# Row 48: Maker
$ws.Range("B48").Value = '''Maker'
$ws.Range("C48").Value = '''https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '''1.036.22'
$ws.Range("E48").Value = '''  +5.03%  '

# Row 49: Aptos
$ws.Range("D49").Value = '''7.683'
$ws.Range("E49").Value = '''  +1.81%  '

# Row 50: EnergySwap
$ws.Range("B50").Value = '''EnergySwap'
$ws.Range("C50").Value = '''https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '''9.971'
$ws.Range("E50").Value = '''  +2.17%  '

# Row 51: RocketPoolETH
$ws.Range("B51").Value = '''RocketPoolETH'
$ws.Range("C51").Value = '''https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '''2.073.03'
$ws.Range("E51").Value = '''  +0.99%  '
